# InitProperty.xlsx maintenance edit:
#   "unify the conception of DataNode, DataTable, Entity."
#
# The sheet that used to represent a single Property/Entity config
# ("Property1") is renamed to the unified "DataNode" naming convention,
# and the in-sheet selection/scroll position left by the author when the
# workbook was last saved is reproduced (frozen header rows 1-8 stay
# frozen, the view is scrolled down and cell E50 is the active selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: Property1 -> DataNode
$ws.Name = "DataNode"

# Reproduce the author's last cursor position/selection in the sheet:
# the frozen pane (header rows 1-8) remains frozen, the window is
# scrolled down so row 24 is at the top of the scrollable area, and the
# active/selected cell is E50.
$excel.ActiveWindow.ScrollRow = 24
$ws.Range("E50").Select()
